$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 49: new journal entry - "Documentation Selection + UML"
$ws.Range("A49").Value = 44693
$ws.Range("B49").Value = 0.33333333333333331
$ws.Range("C49").Value = 0.39930555555555558
$ws.Range("E49").Value = "Documentation Selection + UML"

# Row 50: new journal entry - "Analyse de la partie Placement"
$ws.Range("A50").Value = 44693
$ws.Range("B50").Value = 0.40972222222222227
$ws.Range("C50").Value = 0.51041666666666663
$ws.Range("E50").Value = "Analyse de la partie Placement"

# Row 51: new journal entry - "Explication de chaque phase du placement + algorithme"
$ws.Range("A51").Value = 44693
$ws.Range("B51").Value = 0.5625
$ws.Range("E51").Value = "Explication de chaque phase du placement + algorithme"

# Update the view: scrolled position and active selection
$ws.Application.ActiveWindow.ScrollRow = 46
[void]$ws.Range("F51").Select()
